# Auto-generated: apply scheduled-runner price/profit updates to Marilith_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 275
$ws.Range("I9").Value = 200
$ws.Range("K9").Value = 200
$ws.Range("M9").Value = -31
$ws.Range("H12").Value = 100
$ws.Range("I12").Value = 100
$ws.Range("K12").Value = 100
$ws.Range("M12").Value = 70
$ws.Range("H32").Value = 1041.7142
$ws.Range("J32").Value = 1284.4286
$ws.Range("L32").Value = 1284.4286
$ws.Range("N32").Value = -1936.4286
$ws.Range("H41").Value = 190.25
$ws.Range("I41").Value = 168.3
$ws.Range("K41").Value = 168.3
$ws.Range("M41").Value = 271.7
$ws.Range("H48").Value = 9473
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 9473
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 28419
$ws.Range("M48").Value = $null
$ws.Range("N48").Value = -29003
$ws.Range("H53").Value = 312.72726
$ws.Range("I53").Value = 440.33334
$ws.Range("J53").Value = 159.6
$ws.Range("K53").Value = 440.33334
$ws.Range("L53").Value = 159.6
$ws.Range("M53").Value = 196.66666
$ws.Range("N53").Value = -1433.6
$ws.Range("H56").Value = 9473
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 9473
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 28419
$ws.Range("M56").Value = $null
$ws.Range("N56").Value = -29487
$ws.Range("H62").Value = 4589.8335
$ws.Range("I62").Value = 4271
$ws.Range("J62").Value = 4749.25
$ws.Range("K62").Value = 4271
$ws.Range("L62").Value = 4749.25
$ws.Range("M62").Value = -3647
$ws.Range("N62").Value = -5997.25
$ws.Range("H65").Value = 4589.8335
$ws.Range("I65").Value = 4271
$ws.Range("J65").Value = 4749.25
$ws.Range("K65").Value = 21355
$ws.Range("L65").Value = 23746.25
$ws.Range("M65").Value = -18235
$ws.Range("N65").Value = -29986.25
$ws.Range("H74").Value = 500000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = $null
$ws.Range("H77").Value = 500000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = $null
$ws.Range("H86").Value = 5750
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("H89").Value = 5750
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = $null
$ws.Range("H62").Value = 35000
$ws.Range("J62").Value = 35000
$ws.Range("L62").Value = 35000
$ws.Range("N62").Value = -36248
$ws.Range("H65").Value = 35000
$ws.Range("J65").Value = 35000
$ws.Range("L65").Value = 105000
$ws.Range("N65").Value = -111240
$ws.Range("H132").Value = 1986.3334
$ws.Range("I132").Value = 1986.3334
$ws.Range("K132").Value = 5959.0002
$ws.Range("M132").Value = -3429.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3952.818
$ws.Range("I20").Value = 3648.3
$ws.Range("K20").Value = 3648.3
$ws.Range("M20").Value = -3401.3
$ws.Range("H94").Value = 2338.5
$ws.Range("I94").Value = 2338.5
$ws.Range("K94").Value = 2338.5
$ws.Range("M94").Value = -1887.5
$ws.Range("H134").Value = 4906.2607
$ws.Range("I134").Value = 4792.25
$ws.Range("K134").Value = 14376.75
$ws.Range("M134").Value = -11841.75
$ws.Range("H135").Value = 48000
$ws.Range("J135").Value = 48000
$ws.Range("L135").Value = 48000
$ws.Range("N135").Value = -58140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 450
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -50
$ws.Range("H41").Value = 19796.4
$ws.Range("I41").Value = 3000
$ws.Range("J41").Value = 21662.666
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 21662.666
$ws.Range("M41").Value = -2572
$ws.Range("N41").Value = -22518.666
$ws.Range("H57").Value = 40866.668
$ws.Range("J57").Value = 40866.668
$ws.Range("L57").Value = 40866.668
$ws.Range("N57").Value = -41986.668
$ws.Range("H58").Value = 3102.1177
$ws.Range("I58").Value = 2603.6
$ws.Range("K58").Value = 2603.6
$ws.Range("M58").Value = -2400.6
$ws.Range("H74").Value = 38235.625
$ws.Range("J74").Value = 38235.625
$ws.Range("L74").Value = 38235.625
$ws.Range("N74").Value = -39983.625
$ws.Range("H77").Value = 38235.625
$ws.Range("J77").Value = 38235.625
$ws.Range("L77").Value = 114706.875
$ws.Range("N77").Value = -123442.875
$ws.Range("H107").Value = 948.5
$ws.Range("I107").Value = 899
$ws.Range("K107").Value = 899
$ws.Range("M107").Value = 1021
$ws.Range("H136").Value = 3102.1177
$ws.Range("I136").Value = 2603.6
$ws.Range("K136").Value = 7810.799999999999
$ws.Range("M136").Value = -5260.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 324948.25
$ws.Range("I128").Value = 324948.25
$ws.Range("K128").Value = 974844.75
$ws.Range("M128").Value = -969864.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 11749.75
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 11749.75
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 11749.75
$ws.Range("M33").Value = $null
$ws.Range("N33").Value = -12253.75
$ws.Range("H48").Value = 16245
$ws.Range("I48").Value = 14990
$ws.Range("K48").Value = 14990
$ws.Range("M48").Value = -14505
$ws.Range("H97").Value = 694
$ws.Range("I97").Value = 694
$ws.Range("K97").Value = 694
$ws.Range("M97").Value = -198
$ws.Range("H107").Value = 305.4
$ws.Range("I107").Value = 206
$ws.Range("J107").Value = 703
$ws.Range("K107").Value = 206
$ws.Range("L107").Value = 703
$ws.Range("M107").Value = 1714
$ws.Range("N107").Value = -4543
$ws.Range("H113").Value = 4256.3335
$ws.Range("I113").Value = 2500
$ws.Range("K113").Value = 2500
$ws.Range("M113").Value = -330
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 28800
$ws.Range("J29").Value = 28800
$ws.Range("L29").Value = 28800
$ws.Range("N29").Value = -29390
$ws.Range("H47").Value = 14249.25
$ws.Range("I47").Value = 12000
$ws.Range("K47").Value = 12000
$ws.Range("M47").Value = -11510
$ws.Range("H52").Value = 14249.25
$ws.Range("I52").Value = 12000
$ws.Range("K52").Value = 12000
$ws.Range("M52").Value = -11767
$ws.Range("H61").Value = 3615.8572
$ws.Range("I61").Value = 3906.5
$ws.Range("K61").Value = 3906.5
$ws.Range("M61").Value = -3704.5
$ws.Range("H100").Value = 1625
$ws.Range("I100").Value = 1625
$ws.Range("K100").Value = 1625
$ws.Range("M100").Value = -1084
$ws.Range("H113").Value = 3615.8572
$ws.Range("I113").Value = 3906.5
$ws.Range("K113").Value = 3906.5
$ws.Range("M113").Value = -1736.5
$ws.Range("H122").Value = 3615.375
$ws.Range("J122").Value = 4500
$ws.Range("L122").Value = 13500
$ws.Range("N122").Value = -18400
$ws.Range("H132").Value = 13533.429
$ws.Range("I132").Value = 9530.1875
$ws.Range("K132").Value = 28590.5625
$ws.Range("M132").Value = -26060.5625
$ws.Range("H136").Value = 2104.8262
$ws.Range("I136").Value = 1972.0952
$ws.Range("K136").Value = 5916.2856
$ws.Range("M136").Value = -3366.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 997.5
$ws.Range("I126").Value = 997.5
$ws.Range("K126").Value = 2992.5
$ws.Range("M126").Value = -522.5
$ws.Range("H132").Value = 3487.3333
$ws.Range("I132").Value = 3487.3333
$ws.Range("K132").Value = 10461.9999
$ws.Range("M132").Value = -7931.999899999999
